# Updates cryptos list values (Price / Volume(1h) columns, plus two
# row re-orderings) to match the latest scrape, per commit
# "Updated cryptos list on Wed Mar  6 11:26:02 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.995.59'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '3.878.59'
$ws.Range("E3").Value = '  +4.07%  '
$ws.Range("D4").Value = '''0.997'
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").Value = '''430.50'
$ws.Range("E5").Value = '  +2.34%  '
$ws.Range("D6").Value = '''131.88'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '3.966.82'
$ws.Range("E7").Value = '  +6.64%  '
$ws.Range("D8").Value = '''0.614'
$ws.Range("E8").Value = '  -4.74%  '
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = '''0.731'
$ws.Range("E10").Value = '  -5.24%  '
$ws.Range("D11").Value = '''0.169'
$ws.Range("E11").Value = '  -7.34%  '
$ws.Range("D12").Value = '''0.0000369'
$ws.Range("E12").Value = '  -9.42%  '
$ws.Range("D13").Value = '''41.01'
$ws.Range("E13").Value = '  -4.41%  '
$ws.Range("D14").Value = '4.461.17'
$ws.Range("E14").Value = '  +3.88%  '
$ws.Range("D15").Value = '''10.06'
$ws.Range("E15").Value = '  -4.67%  '
$ws.Range("D16").Value = '''15.68'
$ws.Range("E16").Value = '  +18.45%  '
$ws.Range("D17").Value = '3.845.16'
$ws.Range("E17").Value = '  +3.42%  '
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("D19").Value = '''19.64'
$ws.Range("E19").Value = '  -5.65%  '
$ws.Range("D20").Value = '67.172.57'
$ws.Range("E20").Value = '  +0.64%  '
$ws.Range("D21").Value = '''1.07'
$ws.Range("E21").Value = '  -6.29%  '
$ws.Range("D22").Value = '''409.76'
$ws.Range("E22").Value = '  -7.99%  '
$ws.Range("D23").Value = '''14.49'
$ws.Range("E23").Value = '  -12.26%  '
$ws.Range("D24").Value = '''85.54'
$ws.Range("E24").Value = '  -4.69%  '
$ws.Range("D25").Value = '''3.05'
$ws.Range("E25").Value = '  -3.77%  '
$ws.Range("D26").Value = '''37.61'
$ws.Range("E26").Value = '  -1.98%  '
$ws.Range("D27").Value = '''5.68'
$ws.Range("E27").Value = '  +11.97%  '
$ws.Range("E28").Value = '  -2.38%  '
$ws.Range("D29").Value = '''9.56'
$ws.Range("E29").Value = '  -6.72%  '
$ws.Range("D30").Value = '''690.31'
$ws.Range("E30").Value = '  +4.65%  '
$ws.Range("D31").Value = '''12.47'
$ws.Range("E31").Value = '  -2.94%  '
$ws.Range("D32").Value = '''0.122'
$ws.Range("E32").Value = '  -1.92%  '
$ws.Range("D33").Value = '''2.75'
$ws.Range("E33").Value = '  -0.46%  '
$ws.Range("D34").Value = '''7.17'
$ws.Range("E34").Value = '  -1.60%  '
$ws.Range("D35").Value = '''0.153'
$ws.Range("E35").Value = '  -7.40%  '
$ws.Range("D36").Value = '''38.85'
$ws.Range("E36").Value = '  -7.26%  '
$ws.Range("D37").Value = '0.0₃0813'
$ws.Range("E37").Value = '  +7.68%  '
$ws.Range("D38").Value = '''1.00'
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("D39").Value = '''55.37'
$ws.Range("E39").Value = '  -3.13%  '
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").Value = '''0.0458'
$ws.Range("E41").Value = '  -7.59%  '
$ws.Range("D42").Value = '''0.998'
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = '''0.137'
$ws.Range("E43").Value = '  -8.35%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").Value = '''4.55'
$ws.Range("E44").Value = '  +3.78%  '
$ws.Range("D45").Value = '''148.19'
$ws.Range("E45").Value = '  +0.84%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '''2.09'
$ws.Range("E46").Value = '  -1.89%  '
$ws.Range("B47").Value = 'LidoDAOToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D47").Value = '''3.28'
$ws.Range("E47").Value = '  -5.24%  '
$ws.Range("D48").Value = '''3.11'
$ws.Range("E48").Value = '  -4.03%  '
$ws.Range("D49").Value = '''26.38'
$ws.Range("E49").Value = '  -8.59%  '
$ws.Range("D50").Value = '''2.79'
$ws.Range("E50").Value = '  -3.53%  '
$ws.Range("E51").Value = '  -5.16%  '
